# Rewrite sheet "foo" (sheet1) to reflect the new xlsx2tsv() test fixture:
# adds sep.names-safe header "date.now" / "baz>qux" columns, a percentage
# column, a re-ordered "bar"/colour column, and a right-aligned
# "count of  files" column with some NA / blank entries, plus per-row
# date-format styling (m/d/yyyy, mmm-yy, d-mmm, d-mmm-yy) used to exercise
# detectDates.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("foo")

# --- headers (row 1) ---------------------------------------------------
$ws.Range("A1").Value = "foo"
$ws.Range("B1").Value = "date.now"
$ws.Range("C1").Value = "baz>qux"
$ws.Range("D1").Value = "bar"
$ws.Range("E1").Value = "count of  files"

# --- column A (numeric, unchanged values) ------------------------------
$ws.Range("A2").Value = 170.1
$ws.Range("A3").Value = 10.9
$ws.Range("A4").Value = 199.2
$ws.Range("A5").Value = 229
$ws.Range("A6").Value = 16.4

# --- column B (dates, each row gets a different date format) -----------
$ws.Range("B2").Value = 43532
$ws.Range("B3").Value = 43800
$ws.Range("B4").Value = 43160
$ws.Range("B5").Value = 45337
$ws.Range("B6").Value = 45002

# --- column C (fractions / percentages) ---------------------------------
$ws.Range("C2").Value = 0.28689477956338882
$ws.Range("C3").Value = 0.73915932966927533
$ws.Range("C4").Value = 0.35506826454786744
$ws.Range("C5").Value = 0.98924351240489949
$ws.Range("C6").Value = 0.98532309769208004

# --- column D (colour names, was column B) -------------------------------
$ws.Range("D2").Value = "red"
$ws.Range("D3").Value = "blue"
$ws.Range("D4").Value = "green"
$ws.Range("D5").Value = "black"
$ws.Range("D6").Value = "orange"

# --- column E ("count of  files", right aligned, NA in row 3) ------------
$ws.Range("E2").Value = 16808
$ws.Range("E3").Value = "NA"
$ws.Range("E4").Value = 18630
$ws.Range("E5").Value = ""
$ws.Range("E6").Value = 0

# whole A1:E6 block left-aligned, E column right-aligned (set BEFORE the
# per-cell number formats below so no orphan "number-format only" styles
# get minted ahead of the left-aligned + number-format combo styles)
$ws.Range("A1:E6").HorizontalAlignment = -4131

$ws.Range("B2").NumberFormat = "m/d/yyyy"
$ws.Range("B3").NumberFormat = "mmm-yy"
$ws.Range("B4").NumberFormat = "mmm-yy"
$ws.Range("B5").NumberFormat = "d-mmm"
$ws.Range("B6").NumberFormat = "d-mmm-yy"

$ws.Range("E1:E6").HorizontalAlignment = -4152

$ws.Range("A1").Select()
$ws.Range("C10").Select()

# --- sheet2 -----------------------------------------------------------
$ws2 = $wb.Worksheets.Item("Sheet2")
$ws2.Range("A1").Value = "sheet2"
$ws2.Range("B1").Value = "sheet2_col2"
$ws2.Range("B6").Select()

$ws.Activate()
